$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 2-4 with new values
$ws.Range("A2").Value = "סגן ר' אג""ת"
$ws.Range("B2").Value = "משה"

$ws.Range("A3").Value = "רמ""ח משאבים"
$ws.Range("B3").Value = "אייל"

$ws.Range("A4").Value = "ראש ענף משמעת"
$ws.Range("B4").Value = "ענת"

# Remove rows 5 and 6 (no longer part of the data range)
$ws.Range("A5:B6").Clear()
